# BA_tasks.xlsx edit: add new TODO-list entry (row 13) documenting the
# "tfidf / vocab-size vs. Elasticsearch dense-vector limit" problem, and
# move the view/selection to that area of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- B11 used to borrow its border styling from the slot that we are
# about to repurpose for A13's date formatting below; re-point it at an
# equivalent already-existing "bordered" look (same text, same visual
# appearance) so it keeps rendering the same way.
$ws.Range("B9").Copy($ws.Range("B11"))
$ws.Range("B11").Value = "Matrix aus Bildern (vgl. CIFAR-10), lokale Database mit mehr PDFs füllen"

# --- New row 13: date / task / (notes column skipped) / problem note
$ws.Range("A13").Value = 45151
$ws.Range("A13").NumberFormat = "m/d/yyyy"
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("A13").WrapText = $true

$ws.Range("B13").Value = "tfidf: document term matrix, cosine similarity, tfidf vectorization of ocument corpus"

$ws.Range("D13").Value = "PROBLEM: Elasticsearch max. dimension of dense vector is 2048, vocab size of corpus is bigger (7243)"

$ws.Rows.Item(13).RowHeight = 34

# --- move selection / scroll position to the newly added rows
[void]$ws.Range("C13:C14").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$win.Left = 38420
$win.Top = 840
$win.Width = 34560
$win.Height = 19400
